# Split the bold "July 10 (Wednesday):" run into three runs (the ")"+":"
# piece separated out, and a new "Done" run appended) with a pair of
# <w:proofErr w:type="gramStart"/> / <w:type="gramEnd"/> markers wrapped
# around the "):" + "Done" portion, exactly mirroring the target OOXML
# diff. Word's grammar checker inserts proofErr bookmarks like this when
# text is typed right after a sentence-ending run, which is what produced
# the split here once "Done" was appended after the colon.

$d = $word.ActiveDocument

# Locate the run's text in the body and capture it as a fresh Range so
# InsertXML (which *replaces* a range's contents) lands on exactly the
# right span instead of the collapsed endpoint Find.Execute leaves behind.
$target = $d.Content
$found = $target.Find.Execute("July 10 (Wednesday):", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range = $d.Range($target.Start, $target.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r w:rsidRPr="00535334"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr><w:t>July 10 (Wednesday</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr><w:t>):</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Done</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$range.InsertXML($xml)
